$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.84"
$ws.Range("E2").Value = "'2.19%"
$ws.Range("D3").Value = "'41.44"
$ws.Range("E3").Value = "'2.65%"
$ws.Range("D4").Value = "'5.039"
$ws.Range("E4").Value = "'-0.24%"
$ws.Range("D5").Value = "'0.07548"
$ws.Range("E5").Value = "'3.71%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.598"
$ws.Range("E6").Value = "'2.21%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9270"
$ws.Range("E7").Value = "'0.88%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.421"
$ws.Range("E8").Value = "'3.64%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1218"
$ws.Range("E9").Value = "'5.59%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1838"
$ws.Range("E10").Value = "'6.42%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08970"
$ws.Range("E11").Value = "'3.59%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03943"
$ws.Range("E12").Value = "'-5.63%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1049"
$ws.Range("E13").Value = "'-0.32%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001286"
$ws.Range("E14").Value = "'1.13%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005965"
$ws.Range("E15").Value = "'1.59%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.334"
$ws.Range("E16").Value = "'-1.95%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.385"
$ws.Range("E17").Value = "'2.34%"
$ws.Range("E18").Value = "'1.32%"
$ws.Range("D19").Value = "'7.930"
$ws.Range("E19").Value = "'0.47%"
$ws.Range("D20").Value = "'0.1419"
$ws.Range("E20").Value = "'5.71%"
$ws.Range("D21").Value = "'0.2999"
$ws.Range("E21").Value = "'3.93%"
$ws.Range("D22").Value = "'0.04058"
$ws.Range("E22").Value = "'4.92%"
$ws.Range("D24").Value = "'0.003983"
$ws.Range("E24").Value = "'5.41%"
$ws.Range("E25").Value = "'-4.08%"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("E38").Value = "'4.00%"
$ws.Range("D39").Value = "'0.05214"
$ws.Range("E39").Value = "'5.27%"
$ws.Range("D40").Value = "'0.006400"
$ws.Range("E40").Value = "'-0.49%"
$ws.Range("D41").Value = "'0.007772"
$ws.Range("E41").Value = "'1.21%"
$ws.Range("D42").Value = "'0.1328"
$ws.Range("E42").Value = "'4.46%"
$ws.Range("D43").Value = "'0.007557"
$ws.Range("E43").Value = "'2.31%"
$ws.Range("D44").Value = "'0.007839"
$ws.Range("E44").Value = "'10.80%"
$ws.Range("D45").Value = "'0.3212"
$ws.Range("E45").Value = "'10.76%"
$ws.Range("D46").Value = "'0.00006776"
$ws.Range("E46").Value = "'5.46%"
$ws.Range("E47").Value = "'-0.16%"
$ws.Range("D49").Value = "'0.004200"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("E51").Value = "'-0.16%"
